# Weekly update: insert the newest "Poroto granado" record as a new row
# right above the existing row 120, pushing the rest of the series down
# by one (dimension grows from A1:R131 to A1:R132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 120; existing rows 120-131 shift down
# to 121-132, carrying all their values/formatting with them.
$ws.Rows("120:120").Insert()

# Populate the newly inserted row 120 with the new weekly record.
$ws.Range("A120").Value = 2
$ws.Range("B120").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 45077
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = 100112030
$ws.Range("G120").Value = "Poroto granado"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 400
$ws.Range("K120").Value = 23000
$ws.Range("L120").Value = 25000
$ws.Range("M120").Value = 24000
$ws.Range("N120").Value = "$/malla 25 kilos"
$ws.Range("O120").Value = "Provincia de Limarí"
$ws.Range("P120").Value = 960
$ws.Range("Q120").Value = 25
$ws.Range("R120").Value = "Hortaliza"
